$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.523.79'
$ws.Range('E2').Value = '  -2.13%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.414.19'
$ws.Range('E3').Value = '  -2.79%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.97'
$ws.Range('E5').Value = '  -1.91%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.86'
$ws.Range('E6').Value = '  -2.65%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  -3.06%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.412.67'
$ws.Range('E9').Value = '  -2.88%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.17'
$ws.Range('E10').Value = '  -2.12%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.122'
$ws.Range('E11').Value = '  -3.67%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.435'
$ws.Range('E12').Value = '  -2.75%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.000.89'
$ws.Range('E13').Value = '  -2.77%  '

$ws.Range('E14').Value = '  -0.46%  '

$ws.Range('E15').Value = '  -5.09%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '27.59'
$ws.Range('E16').Value = '  -4.29%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.554.70'
$ws.Range('E17').Value = '  -2.06%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.468.31'
$ws.Range('E18').Value = '  -1.44%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  -2.40%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.79'
$ws.Range('E20').Value = '  -3.98%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '379.09'
$ws.Range('E21').Value = '  -3.52%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.96'
$ws.Range('E22').Value = '  -4.35%  '

$ws.Range('E23').Value = '  -1.83%  '

$ws.Range('E24').Value = '  -0.09%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '71.72'
$ws.Range('E25').Value = '  -2.83%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000118'
$ws.Range('E26').Value = '  -6.43%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.88'
$ws.Range('E27').Value = '  +0.06%  '

$ws.Range('E28').Value = '  -1.47%  '

$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('E30').Value = '  -0.57%  '

$ws.Range('E31').Value = '  -4.15%  '

$ws.Range('E32').Value = '  -3.34%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '23.15'
$ws.Range('E33').Value = '  -2.81%  '

$ws.Range('E34').Value = '  -2.86%  '

$ws.Range('E35').Value = '  -0.13%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '160.84'
$ws.Range('E36').Value = '  -1.50%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.90'
$ws.Range('E37').Value = '  -3.29%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0749'
$ws.Range('E38').Value = '  -3.72%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.877.70'
$ws.Range('E39').Value = '  -7.22%  '

$ws.Range('E40').Value = '  +1.93%  '

$ws.Range('E41').Value = '  -5.14%  '

$ws.Range('E42').Value = '  -0.48%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.52'
$ws.Range('E43').Value = '  -0.72%  '

$ws.Range('E44').Value = '  -3.23%  '

$ws.Range('E45').Value = '  -2.21%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '25.81'
$ws.Range('E46').Value = '  -0.06%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '314.54'
$ws.Range('E47').Value = '  -0.50%  '

$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.22'
$ws.Range('E48').Value = '  -2.25%  '

$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.06'
$ws.Range('E49').Value = '  -5.61%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.49'
$ws.Range('E50').Value = '  -3.55%  '

$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.106'
$ws.Range('E51').Value = '  -3.22%  '
